# Applies the cryptos.xlsx crypto-price/volume update described in the commit
# (GitHub Actions scheduled refresh of coinranking.com data).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.687.64"
$ws.Range("E2").Value = "  -1.27%  "
$ws.Range("D3").Value = "2.228.78"
$ws.Range("E3").Value = "  -0.70%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "252.74"
$ws.Range("E5").Value = "  +8.12%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  -1.22%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "71.18"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.570"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.68"
$ws.Range("E10").Value = "  +18.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0965"
$ws.Range("E11").Value = "  -3.57%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.91"
$ws.Range("E12").Value = "  +0.57%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.106"
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.98"
$ws.Range("E14").Value = "  +1.55%  "
$ws.Range("D15").Value = "2.553.56"
$ws.Range("E15").Value = "  -0.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.96"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.856"
$ws.Range("E17").Value = "  -1.62%  "
$ws.Range("D18").Value = "2.229.57"
$ws.Range("E18").Value = "  -0.60%  "
$ws.Range("D19").Value = "41.599.35"
$ws.Range("E19").Value = "  -1.31%  "
$ws.Range("D20").Value = "0.0₃0969"
$ws.Range("E20").Value = "  -2.32%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.21"
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.04"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.28"
$ws.Range("E23").Value = "  +8.96%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "234.86"
$ws.Range("E24").Value = "  -1.30%  "
$ws.Range("B25").Value = "WEMIXToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.79"
$ws.Range("E25").Value = "  +3.83%  "
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.52"
$ws.Range("E27").Value = "  +6.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.30"
$ws.Range("E28").Value = "  +1.65%  "
$ws.Range("E29").Value = "  +4.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.45"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "20.66"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.121"
$ws.Range("E32").Value = "  +0.70%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.62"
$ws.Range("E33").Value = "  +3.27%  "
$ws.Range("B34").Value = "Stellar"
$ws.Range("C34").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.125"
$ws.Range("E34").Value = "  -2.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0722"
$ws.Range("E35").Value = "  -0.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "26.62"
$ws.Range("E36").Value = "  +18.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.67"
$ws.Range("E37").Value = "  -2.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.00"
$ws.Range("E38").Value = "  +9.99%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0287"
$ws.Range("E39").Value = "  +6.86%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.30"
$ws.Range("E40").Value = "  +1.74%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "69.39"
$ws.Range("E41").Value = "  +1.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.03"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.07"
$ws.Range("E43").Value = "  +18.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.09"
$ws.Range("E44").Value = "  +2.73%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.209"
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.81"
$ws.Range("E46").Value = "  -3.61%  "
$ws.Range("B47").Value = "SynthetixNetwork"
$ws.Range("C47").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.77"
$ws.Range("E47").Value = "  +7.37%  "
$ws.Range("E48").Value = "  +0.73%  "
$ws.Range("E49").Value = "  -0.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.15"
$ws.Range("E50").Value = "  +5.86%  "
$ws.Range("E51").Value = "  +0.73%  "
